$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 = 0, A2 = 0 : bold, centered, top-aligned, thin bordered cells
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0

# B2 = "disconnected_elements" (shared string), default style
$ws.Range("B2").Value = "disconnected_elements"

$r1 = $ws.Range("B1")
$r1.Borders.LineStyle = 1   # xlContinuous
$r1.Borders.Weight = 2      # xlThin
$r1.Font.Bold = $true
$r1.HorizontalAlignment = -4108  # xlCenter
$r1.VerticalAlignment = -4160    # xlTop

# Apply the identical formatting to A2 by copying B1's format so that both
# cells end up sharing the same cell style (matches target workbook).
$r1.Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
